$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new value into D10 (new shared string "websitedev91")
$ws.Range("D10").Value = "websitedev91"

# Set column D width to match diff (stored width 17.88671875 character-units;
# the nearest value this engine's pixel-quantized ColumnWidth setter can produce is 17.833333)
$ws.Columns.Item(4).ColumnWidth = 17

# Update selection to match diff (activeCell D10, sqref D10)
$ws.Range("D10").Select()
